# Update plots for each sample
# The underlying peak-detection threshold for sample S1 / marker CYP2D6_49
# (w_height on peak_table) was lowered from 1000 to 800, which allowed the
# peak to be detected. Propagate the corresponding recalculated results
# into allele_table, marker_table and genotype_result.

$wb = $excel.ActiveWorkbook

# --- peak_table: w_height for CYP2D6_49 (row 4) drops from 1000 to 800 ---
$peakTable = $wb.Worksheets.Item("peak_table")
$peakTable.Range("N4").Value = 800

# --- allele_table: row 6 (CYP2D6_49 / CYP2D6_003) now resolves to a
#     detected, "ok" peak instead of "could not be detected" ---
$alleleTable = $wb.Worksheets.Item("allele_table")
$alleleTable.Range("K6").Value = 800
$alleleTable.Range("L6").Value = 0
$alleleTable.Range("M6").Value = $true
$alleleTable.Range("N6").Value = 20
$alleleTable.Range("O6").Value = 39.3
$alleleTable.Range("P6").Value = 940
$alleleTable.Range("Q6").Value = "ok"
$alleleTable.Range("R6").Value = ""

# --- marker_table: row 4 (CYP2D6_49) now has a called genotype/phenotype ---
$markerTable = $wb.Worksheets.Item("marker_table")
$markerTable.Range("G4").Value = "TT"
$markerTable.Range("H4").Value = "wildtype"

# --- genotype_result: overall sample genotype call ---
$genotypeResult = $wb.Worksheets.Item("genotype_result")
$genotypeResult.Range("B2").Value = "*1/*1"
